# Project Sample Project is saved.TEST Author: admin. Type: SAVE.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# B11 currently holds the text "R40"; change it to the text "1".
# A plain Value/Formula assignment of "1" would be auto-typed as a number,
# so write it as a text formula first and then paste back as a value-only
# (keeps the cell's existing style/number format instead of letting
# NumberFormat mutation fork a new style entry).
$cell = $ws.Range("B11")
$cell.Formula = '="1"'
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false
